$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook is a localization hand-back status report with three sheets:
#   "Overview" - one row per source file, with per-language status columns
#   "zh-cn"    - detailed per-file handoff/handback status for zh-cn
#   "de-de"    - detailed per-file handoff/handback status for de-de
#
# Two files finished processing since the report was last generated:
#   7d77089f-af67-4d13-bf4f-e2576eac4631.md  -> handback succeeded
#   2e08e415-4b1d-4761-aa2a-518002d1e14b.md  -> handback transform failed
# ---------------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"
$statusFailed     = "Handback transform failed"

# --- Overview sheet: row 7 = 7d77089f..., row 8 = 2e08e415... ---------------
$wsOverview.Range("E7").Value = $statusHandedBack
$wsOverview.Range("F7").Value = $statusHandedBack

$wsOverview.Range("E8").Value = $statusFailed
$wsOverview.Range("F8").Value = $statusFailed

# --- zh-cn sheet --------------------------------------------------------
# Row 7 (7d77089f...): handback succeeded
$wsZhCn.Range("C7").Value = $statusHandedBack
$wsZhCn.Range("J7").Value = "7d77089f-af67-4d13-bf4f-e2576eac4631.618a3581a5667c463aae1b1be488427e084e03bc.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-11-29 03:10:24"
$wsZhCn.Range("P7").ClearContents()

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/31da2b384a6fb66389ab4ed9d0b06f193069c858/e2e/7d77089f-af67-4d13-bf4f-e2576eac4631.md",
    [Type]::Missing,
    [Type]::Missing,
    "7d77089f-af67-4d13-bf4f-e2576eac4631.md"
)
$wsZhCn.Range("I7").Style = "HyperLink"

# Row 8 (2e08e415...): handback transform failed
$wsZhCn.Range("C8").Value = $statusFailed
$wsZhCn.Range("P8").Value = "Handback file name: pekldyfg.mzf is different with handoff file name: 2e08e415-4b1d-4761-aa2a-518002d1e14b.31da2b384a6fb66389ab4ed9d0b06f193069c858.zh-cn."

# --- de-de sheet --------------------------------------------------------
# Row 7 (7d77089f...): handback succeeded
$wsDeDe.Range("C7").Value = $statusHandedBack
$wsDeDe.Range("J7").Value = "7d77089f-af67-4d13-bf4f-e2576eac4631.618a3581a5667c463aae1b1be488427e084e03bc.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-11-29 03:10:42"
$wsDeDe.Range("P7").ClearContents()

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/31da2b384a6fb66389ab4ed9d0b06f193069c858/e2e/7d77089f-af67-4d13-bf4f-e2576eac4631.md",
    [Type]::Missing,
    [Type]::Missing,
    "7d77089f-af67-4d13-bf4f-e2576eac4631.md"
)
$wsDeDe.Range("I7").Style = "HyperLink"

# Row 8 (2e08e415...): handback transform failed
$wsDeDe.Range("C8").Value = $statusFailed
$wsDeDe.Range("P8").Value = "Handback file name: pekldyfg.mzf is different with handoff file name: 2e08e415-4b1d-4761-aa2a-518002d1e14b.31da2b384a6fb66389ab4ed9d0b06f193069c858.de-de."

Write-Output "Handback report updated."
